$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "DC Unit Loading Details Name" header cell at F1, matching the
# formatting of the existing header row (row 7).
$ws.Range("A7").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "DC Unit Loading Details Name"

# Add new "Current (DC Units)" data cell at F2, matching the formatting of
# the existing data row (row 8).
$ws.Range("A8").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("F2").Value = "Current (DC Units)"

# Update the active selection to the newly added cells.
$ws.Range("F1:F2").Select()
